# Daily attendance processing - normalize "Recorded By" author ordering.
# For every row in the "Recorded By" column (G) whose value is a
# comma-separated list that currently starts with "system"/"System",
# move the trailing (most-recently-added) author to the front of the list.
# This matches the re-ordering performed by the upstream attendance sync
# job, which always prepends the newest contributor.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$col = 7  # column G = "Recorded By"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $current = $cell.Value2

    if ($current -eq $null) {
        continue
    }

    $parts = $current -split ", "
    $n = $parts.Count

    if ($n -lt 2) {
        continue
    }

    if ($parts[0].ToLower() -ne "system") {
        continue
    }

    $last = $parts[$n - 1]
    $rest = $parts[0..($n - 2)]
    $newParts = @($last) + $rest
    $updated = $newParts -join ", "

    $cell.Value2 = $updated
}

"Recorded By normalization complete"
